# Scheduled runner update: refresh computed Leve profit figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the per-class
# Leve_Profits tables after a fresh market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1254.6875
$ws.Range("I40").Value = 1225.9259
$ws.Range("J40").Value = 1410
$ws.Range("K40").Value = 1225.9259
$ws.Range("L40").Value = 1410
$ws.Range("M40").Value = -1050.9259
$ws.Range("N40").Value = -1760

$ws.Range("H64").Value = 3051.0527
$ws.Range("I64").Value = 2952
$ws.Range("J64").Value = 3161.111
$ws.Range("K64").Value = 2952
$ws.Range("L64").Value = 3161.111
$ws.Range("M64").Value = -2704
$ws.Range("N64").Value = -3657.111

$ws.Range("H67").Value = 3051.0527
$ws.Range("I67").Value = 2952
$ws.Range("J67").Value = 3161.111
$ws.Range("K67").Value = 2952
$ws.Range("L67").Value = 3161.111
$ws.Range("M67").Value = -2094
$ws.Range("N67").Value = -4877.111

$ws.Range("H76").Value = 3936.3157
$ws.Range("I76").Value = 3897
$ws.Range("J76").Value = 4146
$ws.Range("K76").Value = 3897
$ws.Range("L76").Value = 4146
$ws.Range("M76").Value = -3582
$ws.Range("N76").Value = -4776

$ws.Range("H79").Value = 3936.3157
$ws.Range("I79").Value = 3897
$ws.Range("J79").Value = 4146
$ws.Range("K79").Value = 3897
$ws.Range("L79").Value = 4146
$ws.Range("M79").Value = -2805
$ws.Range("N79").Value = -6330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3820
$ws.Range("I63").Value = 3820
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3820
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3134
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 3820
$ws.Range("I66").Value = 3820
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 19100
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -15668
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 2511.5667
$ws.Range("I132").Value = 2013.9565
$ws.Range("J132").Value = 4146.5713
$ws.Range("K132").Value = 6041.8695
$ws.Range("L132").Value = 12439.7139
$ws.Range("M132").Value = -3511.8695
$ws.Range("N132").Value = -17499.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2029.1111
$ws.Range("I105").Value = 2005.7142
$ws.Range("J105").Value = 2111
$ws.Range("K105").Value = 2005.7142
$ws.Range("L105").Value = 2111
$ws.Range("M105").Value = -258.7141999999999
$ws.Range("N105").Value = -5605

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2600
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2600
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2600
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3848

$ws.Range("H65").Value = 2600
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2600
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 13000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -19240

$ws.Range("H132").Value = 1678.826
$ws.Range("I132").Value = 1362.619
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4087.857
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1557.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H57").Value = 20061
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 20061
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 20061
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -21701

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H68").Value = 23863.334

$ws.Range("H70").Value = 4069
$ws.Range("I70").Value = 4106.5264
$ws.Range("J70").Value = 3989.7778
$ws.Range("K70").Value = 4106.5264
$ws.Range("L70").Value = 3989.7778
$ws.Range("M70").Value = -3836.5264
$ws.Range("N70").Value = -4529.7778

$ws.Range("H71").Value = 23863.334

$ws.Range("H73").Value = 4069
$ws.Range("I73").Value = 4106.5264
$ws.Range("J73").Value = 3989.7778
$ws.Range("K73").Value = 4106.5264
$ws.Range("L73").Value = 3989.7778
$ws.Range("M73").Value = -3170.5264
$ws.Range("N73").Value = -5861.7778

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H80").Value = 2499.9375
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 2464.2144
$ws.Range("K80").Value = 2750
$ws.Range("L80").Value = 2464.2144
$ws.Range("M80").Value = -1752
$ws.Range("N80").Value = -4460.2144

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H83").Value = 2499.9375
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 2464.2144
$ws.Range("K83").Value = 13750
$ws.Range("L83").Value = 12321.072
$ws.Range("M83").Value = -8758
$ws.Range("N83").Value = -22305.072

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H94").Value = 49800
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 49800
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 49800
$ws.Range("N94").Value = -51152

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 2590.2222
$ws.Range("I132").Value = 1964.7273
$ws.Range("J132").Value = 3573.1428
$ws.Range("K132").Value = 5894.1819
$ws.Range("L132").Value = 10719.4284
$ws.Range("M132").Value = -3364.1819
$ws.Range("N132").Value = -15779.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3122
$ws.Range("I132").Value = 1373.6833
$ws.Range("J132").Value = 14777.444
$ws.Range("K132").Value = 4121.0499
$ws.Range("L132").Value = 44332.33199999999
$ws.Range("M132").Value = -1591.0499
$ws.Range("N132").Value = -49392.33199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 879.2787
$ws.Range("I132").Value = 742.4151000000001
$ws.Range("J132").Value = 1786
$ws.Range("K132").Value = 2227.2453
$ws.Range("L132").Value = 5358
$ws.Range("M132").Value = 302.7547
$ws.Range("N132").Value = -10418
